# NEL_sitrep_discharge.xlsx update:
#  - appends 20 more days (rows 53-72) of BHRUT/Barts/Homerton discharge
#    sitrep data (with the same running-total formulas used by the rows
#    above them)
#  - adds an XY-scatter chart of the daily combined discharges
#    (DailyBHRUTBartsHom, column H) against the date column

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new daily rows -------------------------------------------------
# Columns: row, Date(A), DailyBHRUT(B), DailyBarts(C), DailyHomerton(D)
# Columns E/F/G/H/I are running totals computed the same way as the
# existing rows (E=CumBHRUT, F=CumBarts, G=CumHomerton,
# H=DailyBHRUTBartsHom, I=CumBHRUTBartsHom).
$rows = @(
  @(53, 43950, 10, 37, 5),
  @(54, 43951, 10, 30, 4),
  @(55, 43952,  7, 32, 1),
  @(56, 43953, 12, 24, 0),
  @(57, 43954,  0, 19, 0),
  @(58, 43955,  9, 28, 3),
  @(59, 43956, 15, 25, 1),
  @(60, 43957, 10, 28, 2),
  @(61, 43958, 12, 23, 1),
  @(62, 43959,  5, 12, 1),
  @(63, 43960,  3, 12, 1),
  @(64, 43961,  2, 13, 0),
  @(65, 43962,  4, 11, 0),
  @(66, 43963,  6, 21, 2),
  @(67, 43964,  6, 13, 3),
  @(68, 43965,  7, 12, 0),
  @(69, 43966,  0, 17, 2),
  @(70, 43967,  2, 11, 2),
  @(71, 43968,  5, 10, 1),
  @(72, 43969,  5, 16, 0)
)

foreach ($row in $rows) {
  $r = $row[0]
  $prev = $r - 1

  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Cells.Item($r, 2).Value = $row[2]
  $ws.Cells.Item($r, 3).Value = $row[3]
  $ws.Cells.Item($r, 4).Value = $row[4]

  $ws.Cells.Item($r, 5).Formula = "=E$prev+B$r"
  $ws.Cells.Item($r, 6).Formula = "=F$prev+C$r"
  $ws.Cells.Item($r, 7).Formula = "=G$prev+D$r"
  $ws.Cells.Item($r, 8).Formula = "=SUM(B" + $r + ":D" + $r + ")"
  $ws.Cells.Item($r, 9).Formula = "=I$prev+H$r"
}

# --- update the view so the newly added rows are visible ------------
[void]$ws.Range("I69").Select()

# --- add an XY scatter chart of DailyBHRUTBartsHom vs Date ----------
$chartObj = $ws.ChartObjects().Add(400, 620, 420, 230)
$chart = $chartObj.Chart
$chart.ChartType = 74  # xlXYScatterLines
$chart.SetSourceData($ws.Range("H1:H67"))

$series = $chart.SeriesCollection(1)
$series.Name = "=Sheet1!`$H`$1"
$series.Formula = "=SERIES(Sheet1!`$H`$1,Sheet1!`$A`$2:`$A`$67,Sheet1!`$H`$2:`$H`$67,1)"

$chart.HasLegend = $false
$chart.HasTitle = $false

Write-Output "Added $($rows.Count) rows and 1 chart"
